# feat: add 2022-Q4 data
#
#  - Insert a new worksheet "2022-Q4" right after "总计" (and before the
#    existing "2022-Q3" sheet), populated with the quarterly fund-holding
#    detail rows.
#  - Update the "总计" (totals) sheet: insert a new leading data row for the
#    2022-Q4 summary and push the older quarters down one row.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$value) {
    # Force a numeric-looking string ("506003", "13.76", ...) to be stored
    # as text instead of being auto-coerced to a number by COM, then strip
    # the number-format style back off so no stray style index is left on
    # the cell (matches cells that were authored as inlineStr with no `s`).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by copying "2022-Q3" (so it inherits
#    identical sheet-level formatting/margins/headers), placed right after
#    "总计", then overwrite its data rows.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $total)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Row 2
$q4.Cells.Item(2, 1).Value = 0
Set-TextValue $q4.Cells.Item(2, 2) "506003"
Set-TextValue $q4.Cells.Item(2, 3) "富国科创板两年定期开放混合"
Set-TextValue $q4.Cells.Item(2, 4) "13.76"
Set-TextValue $q4.Cells.Item(2, 5) "98.91"
Set-TextValue $q4.Cells.Item(2, 6) "2.28"
Set-TextValue $q4.Cells.Item(2, 7) "0.3137"
$q4.Cells.Item(2, 8).Value = 7

# Row 3
$q4.Cells.Item(3, 1).Value = 1
Set-TextValue $q4.Cells.Item(3, 2) "011160"
Set-TextValue $q4.Cells.Item(3, 3) "富国质量成长6个月持有期混合A"
Set-TextValue $q4.Cells.Item(3, 4) "3.97"
Set-TextValue $q4.Cells.Item(3, 5) "89.73"
Set-TextValue $q4.Cells.Item(3, 6) "3.15"
Set-TextValue $q4.Cells.Item(3, 7) "0.1251"
$q4.Cells.Item(3, 8).Value = 4

# Row 4
$q4.Cells.Item(4, 1).Value = 2
Set-TextValue $q4.Cells.Item(4, 2) "004448"
Set-TextValue $q4.Cells.Item(4, 3) "博时汇智回报灵活配置混合"
Set-TextValue $q4.Cells.Item(4, 4) "1.72"
Set-TextValue $q4.Cells.Item(4, 5) "79.42"
Set-TextValue $q4.Cells.Item(4, 6) "3.35"
Set-TextValue $q4.Cells.Item(4, 7) "0.0576"
$q4.Cells.Item(4, 8).Value = 10

# Row 5 (new row - "2022-Q3" only had 3 data rows) - copy the index-column
# ("A") style down from row 4 first so it keeps the shared "s=2" style.
$q4.Range("A4").Copy()
$q4.Range("A5").PasteSpecial(-4122)
$q4.Cells.Item(5, 1).Value = 3
Set-TextValue $q4.Cells.Item(5, 2) "011161"
Set-TextValue $q4.Cells.Item(5, 3) "富国质量成长6个月持有期混合C"
Set-TextValue $q4.Cells.Item(5, 4) "0.15"
Set-TextValue $q4.Cells.Item(5, 5) "89.73"
Set-TextValue $q4.Cells.Item(5, 6) "3.15"
Set-TextValue $q4.Cells.Item(5, 7) "0.0047"
$q4.Cells.Item(5, 8).Value = 4

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing two data rows down one row
#    and write the new 2022-Q4 summary into the freed-up row 2.
# ---------------------------------------------------------------------

# Give row 4 the same index-column ("A") style as the rows above it.
$total.Range("A2").Copy()
$total.Range("A3:A4").PasteSpecial(-4122)

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.1

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.21

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.5
